$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 34.8
$ws.Range("J6").Value = 104.5
$ws.Range("K6").Value = 104.4
$ws.Range("L6").Value = 313.5
$ws.Range("M6").Value = 7.600000000000009
$ws.Range("N6").Value = -537.5
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 1500
$ws.Range("M29").Value = -1219
$ws.Range("H33").Value = 557.4
$ws.Range("I33").Value = 535.61536
$ws.Range("K33").Value = 535.61536
$ws.Range("M33").Value = -306.61536
$ws.Range("H38").Value = 190.16667
$ws.Range("I38").Value = 28.4
$ws.Range("J38").Value = 999
$ws.Range("K38").Value = 85.19999999999999
$ws.Range("L38").Value = 2997
$ws.Range("M38").Value = 286.8
$ws.Range("N38").Value = -3741
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3238
$ws.Range("H60").Value = 1000
$ws.Range("J60").Value = 1000
$ws.Range("L60").Value = 3000
$ws.Range("N60").Value = -3968
$ws.Range("H70").Value = 2877.7778
$ws.Range("I70").Value = 2877.7778
$ws.Range("K70").Value = 8633.3334
$ws.Range("M70").Value = -8363.3334
$ws.Range("H73").Value = 2877.7778
$ws.Range("I73").Value = 2877.7778
$ws.Range("K73").Value = 8633.3334
$ws.Range("M73").Value = -7697.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 478.4
$ws.Range("I2").Value = 573
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 573
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -460
$ws.Range("N2").Value = -326
$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 5000
$ws.Range("K31").Value = 5000
$ws.Range("M31").Value = -4706
$ws.Range("H116").Value = 478.4
$ws.Range("I116").Value = 573
$ws.Range("J116").Value = 100
$ws.Range("K116").Value = 573
$ws.Range("L116").Value = 100
$ws.Range("M116").Value = 1721
$ws.Range("N116").Value = -4688
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 478.4
$ws.Range("I3").Value = 573
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 573
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = -459
$ws.Range("N3").Value = -328
$ws.Range("H32").Value = 40000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""
$ws.Range("H34").Value = 2795
$ws.Range("J34").Value = 2795
$ws.Range("L34").Value = 2795
$ws.Range("N34").Value = -3023
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = ""
$ws.Range("H46").Value = 28332.666
$ws.Range("J46").Value = 32499.5
$ws.Range("L46").Value = 32499.5
$ws.Range("N46").Value = -33095.5
$ws.Range("H94").Value = 2533.3333
$ws.Range("J94").Value = 2800
$ws.Range("L94").Value = 2800
$ws.Range("N94").Value = -3702

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6667700
$ws.Range("I6").Value = 8000240
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 8000240
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = -8000127
$ws.Range("N6").Value = -5226
$ws.Range("H97").Value = 61000
$ws.Range("J97").Value = 61000
$ws.Range("L97").Value = 61000
$ws.Range("N97").Value = -62982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 133.07692
$ws.Range("I23").Value = 41
$ws.Range("J23").Value = 174
$ws.Range("K23").Value = 123
$ws.Range("L23").Value = 522
$ws.Range("M23").Value = 112
$ws.Range("N23").Value = -992
$ws.Range("H33").Value = 64.25
$ws.Range("J33").Value = 96.333336
$ws.Range("L33").Value = 578.000016
$ws.Range("N33").Value = -1144.000016
$ws.Range("H55").Value = 3585.5574
$ws.Range("J55").Value = 3795
$ws.Range("L55").Value = 11385
$ws.Range("N55").Value = -11739
$ws.Range("H92").Value = 344.4
$ws.Range("I92").Value = 344.4
$ws.Range("K92").Value = 1033.2
$ws.Range("M92").Value = 214.8000000000002
$ws.Range("H114").Value = 1813.3334
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 134.11539
$ws.Range("I2").Value = 133.94444
$ws.Range("J2").Value = 134.5
$ws.Range("K2").Value = 133.94444
$ws.Range("L2").Value = 134.5
$ws.Range("M2").Value = -20.94443999999999
$ws.Range("N2").Value = -360.5
$ws.Range("H43").Value = 2800
$ws.Range("I43").Value = 2800
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 2800
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -2649
$ws.Range("N43").Value = ""
$ws.Range("H46").Value = 15699.8
$ws.Range("I46").Value = 8374.75
$ws.Range("J46").Value = 45000
$ws.Range("K46").Value = 8374.75
$ws.Range("L46").Value = 45000
$ws.Range("M46").Value = -8218.75
$ws.Range("N46").Value = -45312
$ws.Range("H98").Value = 8000
$ws.Range("J98").Value = 8000
$ws.Range("L98").Value = 8000
$ws.Range("N98").Value = -13990
$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 4500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -2254
$ws.Range("H122").Value = 5793.3335
$ws.Range("I122").Value = 3977.3333
$ws.Range("J122").Value = 7609.3335
$ws.Range("K122").Value = 11931.9999
$ws.Range("L122").Value = 22828.0005
$ws.Range("M122").Value = -9481.999899999999
$ws.Range("N122").Value = -27728.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1884.25
$ws.Range("I22").Value = 1587
$ws.Range("J22").Value = 1983.3334
$ws.Range("K22").Value = 1587
$ws.Range("L22").Value = 1983.3334
$ws.Range("M22").Value = -1292
$ws.Range("N22").Value = -2573.3334
$ws.Range("H27").Value = 1884.25
$ws.Range("I27").Value = 1587
$ws.Range("J27").Value = 1983.3334
$ws.Range("K27").Value = 1587
$ws.Range("L27").Value = 1983.3334
$ws.Range("M27").Value = -1480
$ws.Range("N27").Value = -2197.3334
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2166.5
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 2249.75
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2249.75
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -4995.75
$ws.Range("H132").Value = 11944.667
$ws.Range("I132").Value = 10474.5
$ws.Range("J132").Value = 14885
$ws.Range("K132").Value = 31423.5
$ws.Range("L132").Value = 44655
$ws.Range("M132").Value = -28893.5
$ws.Range("N132").Value = -49715

